## Adds 4 new PFAS compound rows (PFOA, PFBA, PFBS, PFOS) to Sheet1, mirroring
## the "Add files via upload" commit that appended rows 114-117 below the
## existing 113 rows of data (header in row 1, data rows 2-113).
##
## Columns: A=SMILES, B=MF, C=Name, D:L = 9 molecular descriptors,
## M:V = 10 more molecular descriptors (FNSA-1, MATSv6, ... S7).
## Numeric literals are parsed via [double]"..." to sidestep this host's
## PowerShell-subset lexer, which does not accept scientific-notation
## number literals (e.g. 6E-3) directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function N([string]$s) { return [double]$s }

# Data for the 4 new rows, in column order A..V
$newRows = @(
    @{
        A = "C(=O)(C(C(C(C(C(C(C(F)(F)F)(F)F)(F)F)(F)F)(F)F)(F)F)(F)F)O"
        B = "C8HF15O2"
        C = "PFOA"
        D = N "37.299999999999997"
        E = N "4.4450000000000003"
        F = N "0.89300000000000002"
        G = N "5.5251999999999999"
        H = N "0.53969999999999996"
        I = N "0.53839999999999999"
        J = N "0"
        K = N "0"
        L = N "238.47749887020501"
        M = N "0.88660000000000005"
        N = N "0"
        O = N "6.0000000000000001E-3"
        P = N "28.961300000000001"
        Q = N "4.0174124753451803"
        R = N "-5.6305384885002802E-4"
        S = N "0"
        T = N "-16.865336709401799"
        U = N "20.053894927536199"
        V = N "0"
    },
    @{
        A = "C(=O)(C(C(C(F)(F)F)(F)F)(F)F)O"
        B = "C4HF7O2"
        C = "PFBA"
        D = N "37.299999999999997"
        E = N "1.9039999999999999"
        F = N "9.2999999999999999E-2"
        G = N "2.3693"
        H = N "0.51939999999999997"
        I = N "0.42849999999999999"
        J = N "0"
        K = N "0"
        L = N "0"
        M = N "0.81810000000000005"
        N = N "0"
        O = N "0"
        P = N "25.3216"
        Q = N "12.673150510204"
        R = N "2.7017196235818001E-2"
        S = N "0.26900000000000002"
        T = N "48.988135728567201"
        U = N "22.737689393939299"
        V = N "0"
    },
    @{
        A = "C(C(C(F)(F)S(=O)(=O)O)(F)F)(C(F)(F)F)(F)F"
        B = "C4HF9O3S"
        C = "PFBS"
        D = N "54.37"
        E = N "3.3809999999999998"
        F = N "1.117"
        G = N "3.1642999999999999"
        H = N "0.54290000000000005"
        I = N "0.44169999999999998"
        J = N "0"
        K = N "0"
        L = N "0"
        M = N "0.83479999999999999"
        N = N "0.25800000000000001"
        O = N "0"
        P = N "30.5121"
        Q = N "-1.21474352867955"
        R = N "-5.8348769359496802E-2"
        S = N "-0.04"
        T = N "0.34116008556897198"
        U = N "21.969011350059699"
        V = N "0"
    },
    @{
        A = "C(C(C(C(C(F)(F)S(=O)(=O)O)(F)F)(F)F)(F)F)(C(C(C(F)(F)F)(F)F)(F)F)(F)F"
        B = "C8HF17O3S"
        C = "PFOS"
        D = N "54.37"
        E = N "5.9219999999999997"
        F = N "0.88700000000000001"
        G = N "7.2401999999999997"
        H = N "0.54659999999999997"
        I = N "0.52800000000000002"
        J = N "0"
        K = N "0"
        L = N "233.49042550858999"
        M = N "0.85260000000000002"
        N = N "0"
        O = N "0"
        P = N "51.390300000000003"
        Q = N "3.6760256630087"
        R = N "1.9089778267053E-3"
        S = N "0"
        T = N "-8.9805844080192205"
        U = N "20.059806397306399"
        V = N "0"
    }
)

$colLetters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

$startRow = 114
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    foreach ($col in $colLetters) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
}

# Match the saved view state: scrolled down so row 103 is at top, with the
# (empty) cell below the newly appended data selected.
$excel.ActiveWindow.ScrollRow = 103
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E119").Select()
